$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: B22 was mistakenly entered as text "3" - fix it to the real number 3.
# (The rest of row 22 - A/C/D/E/F/G/H - stays exactly as it was.)
$ws.Cells.Item(22, 2).Value = 3

# Row 23: new annotation row appended right after row 22.
$ws.Cells.Item(23, 1).Value = "Sunsi Wu"

# B23 should hold the numeral "3" stored as text (matching how row 22's
# politeness score used to be entered before the fix above). Excel treats a
# bare numeric-looking string as a number unless the cell is text-formatted,
# so format as text, enter the value, then drop the formatting again so the
# cell keeps its default (un-styled) look while remaining text.
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "3"
$ws.Range("B23").ClearFormats()

$ws.Cells.Item(23, 3).Value = "无"
$ws.Cells.Item(23, 4).Value = "DFT"
$ws.Cells.Item(23, 5).Value = "WRI"
$ws.Cells.Item(23, 6).Value = "01473e7f-4f45-41be-bd3f-03c0ff83190e"
$ws.Cells.Item(23, 7).Value = "H1u8fMW0b_annotated.xlsx"
$ws.Cells.Item(23, 8).Value = "The citations are in non-standard format (section 1.2: Kalman (1960))."
